$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.079.58"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.548.68"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "646.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.394"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").Value = "3.546.19"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "4.230.43"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "94.978.24"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000250"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "3.554.94"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "500.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.466"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000191"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "3.740.73"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.550"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "555.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.892"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +32.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0405"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
